$ppLayoutText = 2
$p = $ppt.ActivePresentation

# Slide 1: Project Initiation & Planning
$s = $p.Slides.Add(1, $ppLayoutText)
[void]$s.Shapes.Item(1).TextFrame.TextRange.InsertAfter(' Project Initiation & Planning ')
[void]$s.Shapes.Item(2).TextFrame.TextRange.InsertAfter("Duration  1 week `nDependencies:  None `nStatus:  Not Started `nResources:  Project Manager ")

# Slide 2: Requirements Gathering & Analysis
$s = $p.Slides.Add(2, $ppLayoutText)
[void]$s.Shapes.Item(1).TextFrame.TextRange.InsertAfter(' Requirements Gathering & Analysis ')
[void]$s.Shapes.Item(2).TextFrame.TextRange.InsertAfter("Duration  2 weeks `nDependencies:  Project Initiation & Planning `nStatus:  Not Started `nResources:  Project Manager, Development Team ")

# Slide 3: System Architecture Design
$s = $p.Slides.Add(3, $ppLayoutText)
[void]$s.Shapes.Item(1).TextFrame.TextRange.InsertAfter(' System Architecture Design ')
[void]$s.Shapes.Item(2).TextFrame.TextRange.InsertAfter("Duration  2 weeks `nDependencies:  Requirements Gathering & Analysis `nStatus:  Not Started `nResources:  Development Team ")

# Slide 4: UI/UX Design
$s = $p.Slides.Add(4, $ppLayoutText)
[void]$s.Shapes.Item(1).TextFrame.TextRange.InsertAfter(' UI/UX Design ')
[void]$s.Shapes.Item(2).TextFrame.TextRange.InsertAfter("Duration  3 weeks `nDependencies:  System Architecture Design `nStatus:  Not Started `nResources:  UI/UX Designer ")

# Slide 5: Backend Development
$s = $p.Slides.Add(5, $ppLayoutText)
[void]$s.Shapes.Item(1).TextFrame.TextRange.InsertAfter(' Backend Development ')
[void]$s.Shapes.Item(2).TextFrame.TextRange.InsertAfter("Duration  12 weeks `nDependencies:  System Architecture Design `nStatus:  Not Started `nResources:  Development Team ")

# Slide 6: iOS App Development
$s = $p.Slides.Add(6, $ppLayoutText)
[void]$s.Shapes.Item(1).TextFrame.TextRange.InsertAfter(' iOS App Development ')
[void]$s.Shapes.Item(2).TextFrame.TextRange.InsertAfter("Duration  10 weeks `nDependencies:  UI/UX Design, Backend Development `nStatus:  Not Started `nResources:  iOS Development Team ")

# Slide 7: Android App Development
$s = $p.Slides.Add(7, $ppLayoutText)
[void]$s.Shapes.Item(1).TextFrame.TextRange.InsertAfter(' Android App Development ')
[void]$s.Shapes.Item(2).TextFrame.TextRange.InsertAfter("Duration  10 weeks `nDependencies:  UI/UX Design, Backend Development `nStatus:  Not Started `nResources:  Android Development Team ")

# Slide 8: Payment Gateway Integration
$s = $p.Slides.Add(8, $ppLayoutText)
[void]$s.Shapes.Item(1).TextFrame.TextRange.InsertAfter(' Payment Gateway Integration ')
[void]$s.Shapes.Item(2).TextFrame.TextRange.InsertAfter("Duration  4 weeks `nDependencies:  Backend Development `nStatus:  Not Started `nResources:  Development Team ")

# Slide 9: Quality Assurance & Testing
$s = $p.Slides.Add(9, $ppLayoutText)
[void]$s.Shapes.Item(1).TextFrame.TextRange.InsertAfter(' Quality Assurance & Testing ')
[void]$s.Shapes.Item(2).TextFrame.TextRange.InsertAfter("Duration  6 weeks `nDependencies:  iOS App Development, Android App Development, Payment Gateway Integration `nStatus:  Not Started `nResources:  QA Team ")

# Slide 10: User Documentation & Onboarding Guides
$s = $p.Slides.Add(10, $ppLayoutText)
[void]$s.Shapes.Item(1).TextFrame.TextRange.InsertAfter(' User Documentation & Onboarding Guides ')
[void]$s.Shapes.Item(2).TextFrame.TextRange.InsertAfter("Duration  2 weeks `nDependencies:  Quality Assurance & Testing `nStatus:  Not Started `nResources:  Technical Writer ")

# Slide 11: Deployment of Backend
$s = $p.Slides.Add(11, $ppLayoutText)
[void]$s.Shapes.Item(1).TextFrame.TextRange.InsertAfter(' Deployment of Backend ')
[void]$s.Shapes.Item(2).TextFrame.TextRange.InsertAfter("Duration  1 week `nDependencies:  Quality Assurance & Testing `nStatus:  Not Started `nResources:  DevOps Team ")

# Slide 12: App Store Submission (iOS & Android)
$s = $p.Slides.Add(12, $ppLayoutText)
[void]$s.Shapes.Item(1).TextFrame.TextRange.InsertAfter(' App Store Submission (iOS & Android) ')
[void]$s.Shapes.Item(2).TextFrame.TextRange.InsertAfter("Duration  1 week `nDependencies:  Quality Assurance & Testing, Deployment of Backend `nStatus:  Not Started `nResources:  Project Manager ")

# Slide 13: Marketing & Launch Preparation
$s = $p.Slides.Add(13, $ppLayoutText)
[void]$s.Shapes.Item(1).TextFrame.TextRange.InsertAfter(' Marketing & Launch Preparation ')
[void]$s.Shapes.Item(2).TextFrame.TextRange.InsertAfter("Duration  4 weeks `nDependencies:  User Documentation & Onboarding Guides, Deployment of Backend, App Store Submission (iOS & Android) `nStatus:  Not Started `nResources:  Marketing Team ")

# Slide 14: Project Launch
$s = $p.Slides.Add(14, $ppLayoutText)
[void]$s.Shapes.Item(1).TextFrame.TextRange.InsertAfter(' Project Launch ')
[void]$s.Shapes.Item(2).TextFrame.TextRange.InsertAfter("Duration  1 day `nDependencies:  Marketing & Launch Preparation `nStatus:  Not Started `nResources:  Project Manager, Marketing Team ")

# Slide 15: Post-Launch Monitoring & Bug Fixes
$s = $p.Slides.Add(15, $ppLayoutText)
[void]$s.Shapes.Item(1).TextFrame.TextRange.InsertAfter(' Post-Launch Monitoring & Bug Fixes ')
[void]$s.Shapes.Item(2).TextFrame.TextRange.InsertAfter("Duration  Ongoing `nDependencies:  Project Launch `nStatus:  Not Started `nResources:  Development Team, QA Team ")
